$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-6 from 45183 to 45184
$ws.Range("C2:C6").Value = 45184
